$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "Pelo q.acompanhei ele fez o serviço mas uma das geladeiras continua sem baixar a temperatura ela chega no mínimo 10 graus "
$ws.Range("C8").Value = 45897.70520351852
$ws.Range("C8").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("D8").Value = "ZjJlMGMyZjMtOWY2Yy00Yzc4LTllZmItMGM5MDE5ZjRjZmU1OjU3MDE2"

$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "'"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = 45897.70435114583
$ws.Range("C9").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("D9").Value = "ZGZmYTBjMjEtN2ZiOC00NTUzLTg2NDEtODg4ZWJlNTI2ZWRjOjU3MDE2"
